$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 50000
$ws.Range("I21").Value = 50000
$ws.Range("K21").Value = 50000
$ws.Range("M21").Value = -49532
# Row 23
$ws.Range("H23").Value = 50000
$ws.Range("I23").Value = 50000
$ws.Range("K23").Value = 50000
$ws.Range("M23").Value = -49766
# Row 31
$ws.Range("H31").Value = 18141.25
$ws.Range("I31").Value = 18141.25
$ws.Range("K31").Value = 54423.75
$ws.Range("M31").Value = -54193.75
# Row 112
$ws.Range("H112").Value = 3006.476
$ws.Range("J112").Value = 3219.2632
$ws.Range("L112").Value = 9657.7896
$ws.Range("N112").Value = -11873.7896
# Row 129
$ws.Range("H129").Value = 1087.6438
$ws.Range("I129").Value = 446.83334
$ws.Range("K129").Value = 1340.50002
$ws.Range("M129").Value = 3659.49998
# Row 133
$ws.Range("H133").Value = 71350
$ws.Range("J133").Value = 71350
$ws.Range("L133").Value = 71350
$ws.Range("N133").Value = -81470

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 132
$ws.Range("H132").Value = 7268.846
$ws.Range("I132").Value = 9492.733
$ws.Range("J132").Value = 4236.273
$ws.Range("K132").Value = 28478.199
$ws.Range("L132").Value = 12708.819
$ws.Range("M132").Value = -25948.199
$ws.Range("N132").Value = -17768.819

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1324337.5
$ws.Range("I58").Value = 1853401.9
$ws.Range("J58").Value = 1676.75
$ws.Range("K58").Value = 1853401.9
$ws.Range("L58").Value = 1676.75
$ws.Range("M58").Value = -1853198.9
$ws.Range("N58").Value = -2082.75
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 134
$ws.Range("H134").Value = 2791.2
$ws.Range("I134").Value = 2239
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6717
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -4182
$ws.Range("N134").Value = -20070
# Row 136
$ws.Range("H136").Value = 1324337.5
$ws.Range("I136").Value = 1853401.9
$ws.Range("J136").Value = 1676.75
$ws.Range("K136").Value = 5560205.699999999
$ws.Range("L136").Value = 5030.25
$ws.Range("M136").Value = -5557655.699999999
$ws.Range("N136").Value = -10130.25

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 243.68182
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 258.05
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 774.1500000000001
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -1244.15
# Row 68
$ws.Range("H68").Value = 1152.5106
$ws.Range("I68").Value = 831.4423
$ws.Range("K68").Value = 2494.3269
$ws.Range("M68").Value = -1683.3269
# Row 71
$ws.Range("H71").Value = 1152.5106
$ws.Range("I71").Value = 831.4423
$ws.Range("K71").Value = 7482.9807
$ws.Range("M71").Value = -3426.9807
# Row 113
$ws.Range("H113").Value = 688.32074
$ws.Range("I113").Value = 585.1212
$ws.Range("J113").Value = 858.6
$ws.Range("K113").Value = 1755.3636
$ws.Range("L113").Value = 2575.8
$ws.Range("M113").Value = 414.6363999999999
$ws.Range("N113").Value = -6915.8
# Row 122
$ws.Range("H122").Value = 609.7895
$ws.Range("I122").Value = 539.13336
$ws.Range("J122").Value = 874.75
$ws.Range("K122").Value = 4852.20024
$ws.Range("L122").Value = 7872.75
$ws.Range("M122").Value = -2402.20024
$ws.Range("N122").Value = -12772.75
# Row 123
$ws.Range("H123").Value = 5100
$ws.Range("I123").Value = 200
$ws.Range("K123").Value = 600
$ws.Range("M123").Value = 1850
# Row 131
$ws.Range("H131").Value = 3018.9038
$ws.Range("J131").Value = 3773.2
$ws.Range("L131").Value = 11319.6
$ws.Range("N131").Value = -21399.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 33819
$ws.Range("J5").Value = 63750
$ws.Range("L5").Value = 63750
$ws.Range("N5").Value = -63974
# Row 24
$ws.Range("H24").Value = 1000000
$ws.Range("J24").Value = 1000000
$ws.Range("L24").Value = 1000000
$ws.Range("N24").Value = -1000346
# Row 29
$ws.Range("H29").Value = 175000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 33
$ws.Range("H33").Value = 17000
$ws.Range("J33").Value = 6000
$ws.Range("L33").Value = 6000
$ws.Range("N33").Value = -6504
# Row 80
$ws.Range("H80").Value = 2759.1304
$ws.Range("J80").Value = 3350
$ws.Range("L80").Value = 3350
$ws.Range("N80").Value = -5346
# Row 83
$ws.Range("H83").Value = 2759.1304
$ws.Range("J83").Value = 3350
$ws.Range("L83").Value = 16750
$ws.Range("N83").Value = -26734

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 25
$ws.Range("H25").Value = 168000
$ws.Range("I25").Value = 168000
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 168000
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -167770
$ws.Range("N25").ClearContents()
# Row 68
$ws.Range("H68").Value = 2173.9546
$ws.Range("I68").Value = 1614
$ws.Range("K68").Value = 1614
$ws.Range("M68").Value = -865
# Row 71
$ws.Range("H71").Value = 2173.9546
$ws.Range("I71").Value = 1614
$ws.Range("K71").Value = 8070
$ws.Range("M71").Value = -4326
# Row 82
$ws.Range("H82").Value = 1713.8636
$ws.Range("I82").Value = 1106
$ws.Range("J82").Value = 3780.6
$ws.Range("K82").Value = 1106
$ws.Range("L82").Value = 3780.6
$ws.Range("M82").Value = -745
$ws.Range("N82").Value = -4502.6
# Row 85
$ws.Range("H85").Value = 1713.8636
$ws.Range("I85").Value = 1106
$ws.Range("J85").Value = 3780.6
$ws.Range("K85").Value = 1106
$ws.Range("L85").Value = 3780.6
$ws.Range("M85").Value = 142
$ws.Range("N85").Value = -6276.6
# Row 132
$ws.Range("H132").Value = 4726.773
$ws.Range("I132").Value = 4829.951
$ws.Range("J132").Value = 3316.6667
$ws.Range("K132").Value = 14489.853
$ws.Range("L132").Value = 9950.000100000001
$ws.Range("M132").Value = -11959.853
$ws.Range("N132").Value = -15010.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4694.6665
$ws.Range("I62").Value = 4620
$ws.Range("J62").Value = 4780
$ws.Range("K62").Value = 4620
$ws.Range("L62").Value = 4780
$ws.Range("M62").Value = -3996
$ws.Range("N62").Value = -6028
# Row 65
$ws.Range("H65").Value = 4694.6665
$ws.Range("I65").Value = 4620
$ws.Range("J65").Value = 4780
$ws.Range("K65").Value = 23100
$ws.Range("L65").Value = 23900
$ws.Range("M65").Value = -19980
$ws.Range("N65").Value = -30140
# Row 68
$ws.Range("H68").Value = 50000
$ws.Range("J68").Value = 50000
$ws.Range("L68").Value = 50000
$ws.Range("N68").Value = -51622
# Row 71
$ws.Range("H71").Value = 50000
$ws.Range("J71").Value = 50000
$ws.Range("L71").Value = 150000
$ws.Range("N71").Value = -158112
# Row 123
$ws.Range("H123").Value = 22784.516
$ws.Range("J123").Value = 22784.516
$ws.Range("L123").Value = 22784.516
$ws.Range("N123").Value = -32584.516
# Row 132
$ws.Range("H132").Value = 2020.7632
$ws.Range("I132").Value = 1997.7084
$ws.Range("J132").Value = 2060.2856
$ws.Range("K132").Value = 5993.1252
$ws.Range("L132").Value = 6180.8568
$ws.Range("M132").Value = -3463.1252
$ws.Range("N132").Value = -11240.8568
# Row 136
$ws.Range("H136").Value = 2642.7646
$ws.Range("I136").Value = 2383.2856
$ws.Range("J136").Value = 3061.923
$ws.Range("K136").Value = 7149.8568
$ws.Range("L136").Value = 9185.769
$ws.Range("M136").Value = -4599.8568
$ws.Range("N136").Value = -14285.769
